$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two row re-orderings)
# to match the latest GitHub Actions scrape.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '93.278.81'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.70%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.419.95'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.26%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '620.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.75%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.38'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.393'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.49%  '

$ws.Range('E9').Value = '  +0.05%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.969'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.29%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.417.17'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.30%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.02'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.53%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.198'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.15%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.22'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.10%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.056.58'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.63%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.052.29'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.37%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000246'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.70%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.18'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.65%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.415.38'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.51%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.15%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.65'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.33%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '498.36'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.08%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.34'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.445'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -9.87%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.60'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.74%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000185'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.60%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '94.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.17%  '

$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.97'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.35%  '

$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.595.44'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.44%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '11.40'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.93%  '

$ws.Range('E31').Value = '  +0.06%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.73'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.15%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.136'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.68%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.02'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.35%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.173'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.80%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '29.99'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.29%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.545'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.51%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '550.16'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.68%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.47'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.18%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.40'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.36%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.922'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.66%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.149'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.40%  '

$ws.Range('B44').Value = 'MantraDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.75'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.71%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.72'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.30%  '

$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.68'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.81%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.51'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.44%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0409'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.80%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.73'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.09%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.09%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.06'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.57%  '
